$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cell contents (keep formatting on header row), then fully clear the rows
# that will no longer be used (7-11) so they disappear from the used range.
$ws.Range("A1:E11").ClearContents()
$ws.Range("A7:E11").Clear()

# --- Header row (row 1) ---
$ws.Cells.Item(1,1).Value = "Identifier"
$ws.Cells.Item(1,2).Value = "Category"
$ws.Cells.Item(1,3).Value = "Rule Number"
$ws.Cells.Item(1,4).Value = "Description"
$ws.Cells.Item(1,5).Value = "Baseline Field Value"
$ws.Cells.Item(1,6).Value = "Candidate Field Value"

# New column F1 needs the same header style as the other header cells.
$ws.Cells.Item(1,5).Copy()
$ws.Cells.Item(1,6).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(1,6).Value = "Candidate Field Value"

$description = "Differences are acceptable if < 3, warnings if 3 < difference < 10, fatal if > 10"

# --- Data rows ---
$ws.Cells.Item(2,1).Value = "BHP       "
$ws.Cells.Item(2,2).Value = "FATAL"
$ws.Cells.Item(2,3).Value = "tolerance_check_1"
$ws.Cells.Item(2,4).Value = $description
$ws.Cells.Item(2,5).Value = 45.2
$ws.Cells.Item(2,6).Value = 4.2

$ws.Cells.Item(3,1).Value = "WBC       "
$ws.Cells.Item(3,2).Value = "FATAL"
$ws.Cells.Item(3,3).Value = "tolerance_check_1"
$ws.Cells.Item(3,4).Value = $description
$ws.Cells.Item(3,5).Value = 24.3
$ws.Cells.Item(3,6).Value = 2

$ws.Cells.Item(4,1).Value = "CBA       "
$ws.Cells.Item(4,2).Value = "FATAL"
$ws.Cells.Item(4,3).Value = "tolerance_check_1"
$ws.Cells.Item(4,4).Value = $description
$ws.Cells.Item(4,5).Value = 103.2
$ws.Cells.Item(4,6).Value = 13.2

$ws.Cells.Item(5,1).Value = "WBC       "
$ws.Cells.Item(5,2).Value = "FATAL"
$ws.Cells.Item(5,3).Value = "tolerance_check_1"
$ws.Cells.Item(5,4).Value = $description
$ws.Cells.Item(5,5).Value = 24.1
$ws.Cells.Item(5,6).Value = 4.1

$ws.Cells.Item(6,1).Value = "ANZ       "
$ws.Cells.Item(6,2).Value = "FATAL"
$ws.Cells.Item(6,3).Value = "tolerance_check_1"
$ws.Cells.Item(6,4).Value = $description
$ws.Cells.Item(6,5).Value = 29.15
$ws.Cells.Item(6,6).Value = 2.15

$ws.Range("A1").Select()
